$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Low): remove B2 and E2, update C2/D2/F2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 2033
$ws.Range("D2").Value = 1886
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 4222

# Row 3 (mid-Low): update all values
$ws.Range("B3").Value = 5259
$ws.Range("C3").Value = 1009.25
$ws.Range("D3").Value = 96.89
$ws.Range("E3").Value = 283
$ws.Range("F3").Value = 5919

# Row 4 (Middle): update all values
$ws.Range("B4").Value = 5559
$ws.Range("C4").Value = 2934.69
$ws.Range("D4").Value = 2283.14
$ws.Range("E4").Value = 1797
$ws.Range("F4").Value = 56.18

# Row 5 (mid-High): remove B5 and E5, update C5/D5/F5
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = 348.06
$ws.Range("D5").Value = 1346.96
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 4509.38

# Row 6 (High): update B6, E6, F6 (C6/D6 remain empty)
$ws.Range("B6").Value = 4476
$ws.Range("E6").Value = 869
$ws.Range("F6").Value = 18.77

Write-Host "Edit applied successfully"
